# "change the order of experience"
# Swap the first two entries on the "experience" sheet (rows 2 and 3)
# so the earlier "Summer Intern" (Jun. 2016) appears before the
# "Ph.D. Candidate" (Sep. 2016 - Jul. 2019) entry, then make the
# experience sheet the active tab/selection and set its page setup.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("experience")

# --- Swap the contents of rows 2 and 3 (columns A:F) ---
$columns = @("A", "B", "C", "D", "E", "F")
foreach ($col in $columns) {
    $cellRow2 = $col + "2"
    $cellRow3 = $col + "3"
    $valueRow2 = $ws.Range($cellRow2).Value2
    $valueRow3 = $ws.Range($cellRow3).Value2
    $ws.Range($cellRow2).Value2 = $valueRow3
    $ws.Range($cellRow3).Value2 = $valueRow2
}

# --- Swap the row heights to match the (now swapped) content ---
$heightRow2 = $ws.Rows.Item(2).RowHeight
$heightRow3 = $ws.Rows.Item(3).RowHeight
$ws.Rows.Item(2).RowHeight = $heightRow3
$ws.Rows.Item(3).RowHeight = $heightRow2

# --- Page setup for the experience sheet (A4, portrait) ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Make "experience" the active sheet/tab with a new selection ---
$ws.Activate() | Out-Null
$ws.Range("A5").Select() | Out-Null
